# Rectify the sample ID columns: cells in columns B and C of these rows
# currently hold the Sample ID as a text string with a trailing "R"
# (e.g. "1022953R"). Replace them with the plain numeric sample id
# (e.g. 1022953) in both columns, turning the cell from a text/shared
# -string cell into a genuine number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    80  = 1022953
    88  = 1022930
    92  = 1022935
    96  = 1022944
    103 = 1022941
    108 = 1022939
    111 = 1022942
    112 = 1022933
    122 = 1022938
    130 = 1022940
    131 = 1022927
    134 = 1022928
    137 = 1022932
    160 = 1022949
    165 = 1022934
    166 = 1022946
    174 = 1022929
    176 = 1022937
    196 = 1022943
    205 = 1022945
    210 = 1022948
}

foreach ($r in $rows.Keys) {
    $id = $rows[$r]
    $ws.Range("B$r").Value = $id
    $ws.Range("C$r").Value = $id
}

# Leave the cursor where the last edit was made (matches the saved view
# state after this edit: active cell D210).
[void]$ws.Range("D210").Select()
